$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct the attribute value in B4 (was "PSet_CN", should be "PSet_PN")
$ws.Range("B4").Value = "PSet_PN"

# Update the view: zoom in further and move selection to C6, scrolled to top
$excel.ActiveWindow.Zoom = 280
$ws.Range("C6").Select()
